$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update date column (A2:A5) from 2025-11-29 to 2025-12-01.
# Temporarily force text format so Excel keeps the value as plain text
# (shared string) rather than auto-converting it to a date serial number,
# then restore the original (unstyled) cell style so no formatting diff
# is introduced.
$dateRange = $ws.Range("A2:A5")
$dateRange.NumberFormat = "@"
$dateRange.Value = "2025-12-01"
$dateRange.Style = "Normal"

# Update MACRO_SCORE column (N2:N5)
$ws.Range("N2:N5").Value = 85.87246918135976
